# feat: add 2022-Q1 data
#
# - Insert a new "2022-Q1" worksheet (same layout as the other quarterly
#   sheets) positioned right after "2021-Q2" and before "总计".
# - Update the "总计" (grand total) sheet with a new leading row for the
#   2022-Q1 quarter, pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet and position it before "总计".
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q1"

$total = $wb.Worksheets.Item("总计")
$newSheet.Move($total)

# Re-fetch by name: after Move() the old handle can go stale.
$q1 = $wb.Worksheets.Item("2022-Q1")

# Use "2021-Q2" as the style template - it already carries the header /
# index-column formatting ("s=2": bold, centered, bordered) that the new
# sheet should copy.
$template = $wb.Worksheets.Item("2021-Q2")

$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$template.Range("A2:A3").Copy()
$q1.Range("A2:A3").PasteSpecial(-4122)

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2 - 010591 / 富国中国中小盘混合(QDII)美元
# Fund code / ratios look like numbers, but must be stored as plain text
# (matching the source data) - write with a leading "'" so the host keeps
# them as text, then ClearFormats() to drop the cosmetic "quote prefix"
# style the leading "'" leaves behind (data keeps its text type/value).
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'010591"
$q1.Range("C2").Value = "富国中国中小盘混合(QDII)美元"
$q1.Range("D2").Value = "'35.75"
$q1.Range("E2").Value = "'86.53"
$q1.Range("F2").Value = "'2.29"
$q1.Range("G2").Value = "'0.8187"
$q1.Range("H2").Value = 10
$q1.Range("B2").ClearFormats()
$q1.Range("D2:G2").ClearFormats()

# Row 3 - 100061 / 富国中国中小盘混合(QDII)人民币
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'100061"
$q1.Range("C3").Value = "富国中国中小盘混合(QDII)人民币"
$q1.Range("D3").Value = "'35.75"
$q1.Range("E3").Value = "'86.53"
$q1.Range("F3").Value = "'2.29"
$q1.Range("G3").Value = "'0.8187"
$q1.Range("H3").Value = 10
$q1.Range("B3").ClearFormats()
$q1.Range("D3:G3").ClearFormats()

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: push rows 2-4 down to 3-5, and write the
#    new 2022-Q1 row at row 2.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Snapshot the existing 3 data rows (bottom-to-top) before overwriting.
# NOTE: `.Value` getter is not wired up in this host - use `.Value2` instead.
$row4 = @($total.Range("A4").Value2, $total.Range("B4").Value2, $total.Range("C4").Value2, $total.Range("D4").Value2)
$row3 = @($total.Range("A3").Value2, $total.Range("B3").Value2, $total.Range("C3").Value2, $total.Range("D3").Value2)
$row2 = @($total.Range("A2").Value2, $total.Range("B2").Value2, $total.Range("C2").Value2, $total.Range("D2").Value2)

# Copy the formatting of the last existing data row down into the new row 5.
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)

# Shift old row 3 -> row 4, old row 2 -> row 3 (values only, keep formats).
$total.Range("A5").Value = 3
$total.Range("B5").Value = $row4[1]
$total.Range("C5").Value = $row4[2]
$total.Range("D5").Value = $row4[3]

$total.Range("A4").Value = 2
$total.Range("B4").Value = $row3[1]
$total.Range("C4").Value = $row3[2]
$total.Range("D4").Value = $row3[3]

$total.Range("A3").Value = 1
$total.Range("B3").Value = $row2[1]
$total.Range("C3").Value = $row2[2]
$total.Range("D3").Value = $row2[3]

# New row 2: 2022-Q1
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 1.64

Write-Output "2022-Q1 sheet added; 总计 sheet updated"
